# This script reproduces the target edit of the "Artfynd" worksheet:
#  - The existing 12 data rows (rows 2-13) are reordered into a new layout
#    (rows 2-7, 11-16), with their full content (all columns) carried along
#    unchanged.
#  - Three brand-new observation rows are inserted at rows 8-10.
#  - The sheet dimension grows from A1:AY13 to A1:AY16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of the original row number (in the workbook as it was opened)
# to its new row number after the reorder.
$beforeToAfter = [ordered]@{
    2  = 11
    3  = 2
    4  = 3
    5  = 12
    6  = 4
    7  = 5
    8  = 6
    9  = 7
    10 = 13
    11 = 14
    12 = 15
    13 = 16
}

$lastCol = "AY"
$stageOffset = 1000

# Phase 1: move every existing row out of the way into a staging area
# (original row number + 1000) so that the forward/backward moves below
# never collide with each other while they are in flight.
foreach ($b in $beforeToAfter.Keys) {
    $stageRow = $b + $stageOffset
    $src = $ws.Range("A" + $b + ":" + $lastCol + $b)
    $dst = $ws.Range("A" + $stageRow)
    $src.Cut($dst)
}

# Phase 2: move every row from its staging position to its final position.
foreach ($b in $beforeToAfter.Keys) {
    $stageRow = $b + $stageOffset
    $finalRow = $beforeToAfter[$b]
    $src = $ws.Range("A" + $stageRow + ":" + $lastCol + $stageRow)
    $dst = $ws.Range("A" + $finalRow)
    $src.Cut($dst)
}

# Phase 3: populate the three new rows (8, 9, 10) with the newly reported
# observations.
$newRows = @(
    @{
        Row = 8
        A = 112379079; B = 90800; C = "Ovaliderad"; D = "LC"; E = 4364
        F = "Dropptaggsvamp"; G = "Hydnellum ferrugineum"; H = "(Fr.:Fr.) P. Karst."
        P = "Lortmossen, Vstm"; Q = 532226; R = 6623334; S = 25
        T = "Västmanland"; U = "Skinnskatteberg"; V = "Västmanland"; W = "Skinnskatteberg"
        Y = "2023-08-27"; AA = "2023-09-18"
        AW = "Mikael Hagström"; AX = "Mikael Hagström"
    },
    @{
        Row = 9
        A = 112379172; B = 90800; C = "Ovaliderad"; D = "LC"; E = 4364
        F = "Dropptaggsvamp"; G = "Hydnellum ferrugineum"; H = "(Fr.:Fr.) P. Karst."
        P = "Lortmossen, Vstm"; Q = 532294; R = 6623516; S = 25
        T = "Västmanland"; U = "Skinnskatteberg"; V = "Västmanland"; W = "Skinnskatteberg"
        Y = "2023-08-27"; AA = "2023-09-18"
        AW = "Mikael Hagström"; AX = "Mikael Hagström"
    },
    @{
        Row = 10
        A = 112379171; B = 90800; C = "Ovaliderad"; D = "LC"; E = 4364
        F = "Dropptaggsvamp"; G = "Hydnellum ferrugineum"; H = "(Fr.:Fr.) P. Karst."
        P = "Lortmossen, Vstm"; Q = 532169; R = 6623457; S = 25
        T = "Västmanland"; U = "Skinnskatteberg"; V = "Västmanland"; W = "Skinnskatteberg"
        Y = "2023-08-27"; AA = "2023-09-18"
        AW = "Mikael Hagström"; AX = "Mikael Hagström"
    }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row
    $ws.Range("A" + $r).Value = $rowData.A
    $ws.Range("B" + $r).Value = $rowData.B
    $ws.Range("C" + $r).Value = $rowData.C
    $ws.Range("D" + $r).Value = $rowData.D
    $ws.Range("E" + $r).Value = $rowData.E
    $ws.Range("F" + $r).Value = $rowData.F
    $ws.Range("G" + $r).Value = $rowData.G
    $ws.Range("H" + $r).Value = $rowData.H
    $ws.Range("P" + $r).Value = $rowData.P
    $ws.Range("Q" + $r).Value = $rowData.Q
    $ws.Range("R" + $r).Value = $rowData.R
    $ws.Range("S" + $r).Value = $rowData.S
    $ws.Range("T" + $r).Value = $rowData.T
    $ws.Range("U" + $r).Value = $rowData.U
    $ws.Range("V" + $r).Value = $rowData.V
    $ws.Range("W" + $r).Value = $rowData.W
    $ws.Range("Y" + $r).Value = $rowData.Y
    $ws.Range("AA" + $r).Value = $rowData.AA
    $ws.Range("AD" + $r).Value = $false
    $ws.Range("AE" + $r).Value = $false
    $ws.Range("AG" + $r).Value = $false
    $ws.Range("AW" + $r).Value = $rowData.AW
    $ws.Range("AX" + $r).Value = $rowData.AX
}
